$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @("0.1184146901506048", "0.5185973060008381", "0.4691159408045538", "0.6849203901217672", "0.698284046247675", "15")
    3  = @("0.1872890445127373", "0.5754152234511226", "0.6110775347743189", "0.7817144841784108", "0.7875963862927191", "14")
    4  = @("0.1843166320775897", "0.6273081474520239", "0.5779702429309931", "0.7602435418541831", "0.7676788112505872", "13")
    5  = @("0.3221619536241352", "0.5655998917628668", "0.5382575468002283", "0.7336603756509059", "0.6884521947652315", "12")
    6  = @("0.2262967854219969", "0.4712852108734478", "0.2956388450555676", "0.5437268110508876", "0.5185281775945857", "11")
    7  = @("0.3318332233303713", "0.399431789294712", "0.2898093932768167", "0.53833947772462", "0.4468359196381605", "10")
    8  = @("0.3032448163457492", "0.4442789123723634", "0.3348453944406773", "0.5786582708651776", "0.5227322189947009", "9")
    9  = @("0.3395226469355723", "0.4954646948149822", "0.3955394669362153", "0.6289192849135852", "0.5659518799878692", "8")
    10 = @("0.2802560178301312", "0.3574868270362884", "0.1586788892857869", "0.3983451886062977", "0.3057635732311509", "7")
    11 = @("0.3351585394422735", "0.4364782826832794", "0.3756151435636499", "0.612874492505317", "0.5620860044540215", "6")
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = [double]$vals[0]
    $ws.Cells.Item($row, 3).Value = [double]$vals[1]
    $ws.Cells.Item($row, 4).Value = [double]$vals[2]
    $ws.Cells.Item($row, 5).Value = [double]$vals[3]
    $ws.Cells.Item($row, 6).Value = [double]$vals[4]
    $ws.Cells.Item($row, 7).Value = [double]$vals[5]
}
